$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 with "Save" label, using same style as other header cells (copy from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for rows 2-8
$values = @(0, 1, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
